$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44627
$ws.Range("J2").Value2 = 120
$ws.Range("K2").Value2 = 4000
$ws.Range("L2").Value2 = 4500
$ws.Range("M2").Value2 = 4250
$ws.Range("P2").Value2 = 71

$ws.Range("D3").Value2 = 44382
$ws.Range("J3").Value2 = 160
$ws.Range("K3").Value2 = 7000
$ws.Range("L3").Value2 = 8000
$ws.Range("M3").Value2 = 7438
$ws.Range("P3").Value2 = 124

$ws.Range("D4").Value2 = 44494
$ws.Range("J4").Value2 = 120
$ws.Range("K4").Value2 = 5000
$ws.Range("L4").Value2 = 6000
$ws.Range("M4").Value2 = 5500
$ws.Range("P4").Value2 = 92

$ws.Range("D5").Value2 = 44740
$ws.Range("J5").Value2 = 120
$ws.Range("K5").Value2 = 6000
$ws.Range("L5").Value2 = 7000
$ws.Range("M5").Value2 = 6500
$ws.Range("P5").Value2 = 108

$ws.Range("D6").Value2 = 44648
$ws.Range("J6").Value2 = 120
$ws.Range("K6").Value2 = 6500
$ws.Range("L6").Value2 = 7000
$ws.Range("M6").Value2 = 6750
$ws.Range("P6").Value2 = 112

$ws.Range("D7").Value2 = 44242
$ws.Range("J7").Value2 = 160
$ws.Range("K7").Value2 = 5000
$ws.Range("L7").Value2 = 5500
$ws.Range("M7").Value2 = 5250
$ws.Range("P7").Value2 = 88

$ws.Range("D8").Value2 = 44669
$ws.Range("J8").Value2 = 130
$ws.Range("K8").Value2 = 4500
$ws.Range("L8").Value2 = 5000
$ws.Range("M8").Value2 = 4750
$ws.Range("P8").Value2 = 79

$ws.Range("D9").Value2 = 44603
$ws.Range("J9").Value2 = 140
$ws.Range("K9").Value2 = 5500
$ws.Range("L9").Value2 = 6000
$ws.Range("M9").Value2 = 5750
$ws.Range("P9").Value2 = 96

$ws.Range("D10").Value2 = 44657
$ws.Range("J10").Value2 = 100
$ws.Range("K10").Value2 = 5000
$ws.Range("L10").Value2 = 5500
$ws.Range("M10").Value2 = 5250
$ws.Range("P10").Value2 = 88

$ws.Range("D11").Value2 = 44676
$ws.Range("J11").Value2 = 120
$ws.Range("K11").Value2 = 4000
$ws.Range("L11").Value2 = 4500
$ws.Range("M11").Value2 = 4250
$ws.Range("P11").Value2 = 71

$ws.Range("D12").Value2 = 44400
$ws.Range("J12").Value2 = 120
$ws.Range("K12").Value2 = 9000
$ws.Range("L12").Value2 = 10000
$ws.Range("M12").Value2 = 9500
$ws.Range("P12").Value2 = 158

$ws.Range("D13").Value2 = 44827
$ws.Range("J13").Value2 = 120
$ws.Range("K13").Value2 = 6000
$ws.Range("L13").Value2 = 7000
$ws.Range("M13").Value2 = 6500
$ws.Range("P13").Value2 = 108

$ws.Range("D14").Value2 = 44760
$ws.Range("J14").Value2 = 130
$ws.Range("K14").Value2 = 7000
$ws.Range("L14").Value2 = 7500
$ws.Range("M14").Value2 = 7250
$ws.Range("P14").Value2 = 121

$ws.Range("D15").Value2 = 44764
$ws.Range("J15").Value2 = 120
$ws.Range("K15").Value2 = 7000
$ws.Range("L15").Value2 = 8000
$ws.Range("M15").Value2 = 7500
$ws.Range("P15").Value2 = 125

$ws.Range("D16").Value2 = 44785
$ws.Range("J16").Value2 = 130
$ws.Range("K16").Value2 = 7000
$ws.Range("L16").Value2 = 8000
$ws.Range("M16").Value2 = 7500
$ws.Range("P16").Value2 = 125

$ws.Range("D17").Value2 = 44589
$ws.Range("J17").Value2 = 110
$ws.Range("K17").Value2 = 5000
$ws.Range("L17").Value2 = 6000
$ws.Range("M17").Value2 = 5500
$ws.Range("P17").Value2 = 92

$ws.Range("D18").Value2 = 44421
$ws.Range("J18").Value2 = 100
$ws.Range("K18").Value2 = 8000
$ws.Range("L18").Value2 = 9000
$ws.Range("M18").Value2 = 8500
$ws.Range("P18").Value2 = 142

$ws.Range("D19").Value2 = 44362
$ws.Range("J19").Value2 = 120
$ws.Range("K19").Value2 = 8000
$ws.Range("L19").Value2 = 9000
$ws.Range("M19").Value2 = 8500
$ws.Range("P19").Value2 = 142

$ws.Range("D20").Value2 = 44281
$ws.Range("J20").Value2 = 120
$ws.Range("K20").Value2 = 5500
$ws.Range("L20").Value2 = 6000
$ws.Range("M20").Value2 = 5750
$ws.Range("P20").Value2 = 96
